# Adiciona politica de preco: insere as colunas "modelo" e "politica"
# entre "full" e as demais, atualiza o link (tracking_id) e os dados
# de tipo/modelo/politica de cada linha.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere duas novas colunas antes da antiga coluna D ("tipo"), empurrando
# full/tipo/link de C/D/E para E/F/G.
$ws.Columns("C:D").Insert()

# Cabecalho (linha 1)
$ws.Range("C1").Value = "modelo"
$ws.Range("D1").Value = "politica"

# Linha 2
$ws.Range("C2").Value = "Sem Modelo"
$ws.Range("D2").Value = ""
$ws.Range("F2").Value = "classico"
$ws.Range("G2").Value = "https://www.mercadolivre.com.br/controle-longa-distncia-jfa-acqua-1200-resistente-a-agua/p/MLB27687422?pdp_filters=seller_id:154235731#searchVariation=MLB27687422&position=3&search_layout=stack&type=product&tracking_id=c7de5016-8467-4dd8-8318-9a4483899c31"

# Linha 3
$ws.Range("C3").Value = "Sem Modelo"
$ws.Range("D3").Value = ""
$ws.Range("F3").Value = "classico"
$ws.Range("G3").Value = "https://www.mercadolivre.com.br/conversor-fio-para-rca-remoto-slim-12v-jfa-automotivo-cd-dvd/p/MLB25707531?pdp_filters=seller_id:154235731#searchVariation=MLB25707531&position=4&search_layout=stack&type=product&tracking_id=c7de5016-8467-4dd8-8318-9a4483899c31"

# Linha 4
$ws.Range("C4").Value = "Sem Modelo"
$ws.Range("D4").Value = ""
$ws.Range("F4").Value = "classico"
$ws.Range("G4").Value = "https://www.mercadolivre.com.br/controle-remoto-universal-longa-distncia-jfa-k1200-preto/p/MLB28687615?pdp_filters=seller_id:154235731#searchVariation=MLB28687615&position=1&search_layout=stack&type=product&tracking_id=c7de5016-8467-4dd8-8318-9a4483899c31"

# Linha 5
$ws.Range("C5").Value = "Sem Modelo"
$ws.Range("D5").Value = ""
$ws.Range("F5").Value = "classico"
$ws.Range("G5").Value = "https://produto.mercadolivre.com.br/MLB-4531110844-filtro-anti-ruido-jfa-com-blindagem-eletromagnetica-2020k-_JM#position%3D5%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dc7de5016-8467-4dd8-8318-9a4483899c31"

# Linha 6
$ws.Range("C6").Value = "Sem Modelo"
$ws.Range("D6").Value = ""
$ws.Range("F6").Value = "classico"
$ws.Range("G6").Value = "https://produto.mercadolivre.com.br/MLB-4531096344-voltimetro-jfa-vs5hi-3-em-1-sequenciador-high-voltagem-12v-_JM#position%3D6%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dc7de5016-8467-4dd8-8318-9a4483899c31"

# Linha 7
$ws.Range("C7").Value = "FONTE 70A"
$ws.Range("D7").Value = "Acima"
$ws.Range("F7").Value = "classico"
$ws.Range("G7").Value = "https://produto.mercadolivre.com.br/MLB-3629883283-fonte-digital-jfa-70a-storm-carregador-inteligente-bateria-_JM#position%3D7%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dc7de5016-8467-4dd8-8318-9a4483899c31"

# Linha 8
$ws.Range("C8").Value = "FONTE 90 BOB"
$ws.Range("D8").Value = "Acima"
$ws.Range("F8").Value = "classico"
$ws.Range("G8").Value = "https://produto.mercadolivre.com.br/MLB-3629903553-fonte-automotiva-jfa-bob-storm-90a-bivolt-carregador-_JM#position%3D8%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dc7de5016-8467-4dd8-8318-9a4483899c31"

# Linha 9
$ws.Range("C9").Value = "FONTE 40A"
$ws.Range("D9").Value = "Acima"
$ws.Range("F9").Value = "classico"
$ws.Range("G9").Value = "https://produto.mercadolivre.com.br/MLB-3629872501-fonte-carregador-de-bateria-jfa-40a-storm-som-automotivo-_JM#position%3D9%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dc7de5016-8467-4dd8-8318-9a4483899c31"

# Linha 10
$ws.Range("C10").Value = "FONTE 60A"
$ws.Range("D10").Value = "Acima"
$ws.Range("F10").Value = "classico"
$ws.Range("G10").Value = "https://produto.mercadolivre.com.br/MLB-3629847295-fonte-carregador-automotivo-jfa-storm-60a-bivolt-voltamp-_JM#position%3D10%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Dc7de5016-8467-4dd8-8318-9a4483899c31"
